$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 531 (shifts existing rows 531..605 down to 532..606)
$ws.Rows.Item(531).Insert()

# Populate the newly inserted row 531 with the new weekly data point
$ws.Cells.Item(531, 1).Value = 4
$ws.Cells.Item(531, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(531, 3).Value = "Los Lagos"
$ws.Cells.Item(531, 4).Value = 45127
$ws.Cells.Item(531, 5).Value = 10
$ws.Cells.Item(531, 6).Value = 100114013
$ws.Cells.Item(531, 7).Value = "Zanahoria"
$ws.Cells.Item(531, 8).Value = "Sin especificar"
$ws.Cells.Item(531, 9).Value = "Primera"
$ws.Cells.Item(531, 10).Value = 250
$ws.Cells.Item(531, 11).Value = 7500
$ws.Cells.Item(531, 12).Value = 8000
$ws.Cells.Item(531, 13).Value = 7700
$ws.Cells.Item(531, 14).Value = '$/saco 20 kilos'
$ws.Cells.Item(531, 15).Value = "Provincia de Llanquihue"
$ws.Cells.Item(531, 16).Value = 385
$ws.Cells.Item(531, 17).Value = 20
$ws.Cells.Item(531, 18).Value = "Hortaliza"

# Match the date-style formatting used by the rest of column D
$ws.Cells.Item(531, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
